# Update the test-data e-mail addresses so their embedded timestamp
# changes from 20251109_004215 to 20251109_005042.
# The addresses live in the "E-Mail" column (C) of the "UsuariosRegistro"
# sheet, rows 2-6. The "LoginData" sheet references the same shared
# strings for its first two rows, so it will pick up the change too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$ws.Range("C2").Value = "juan.perez+20251109_005042@test.com"
$ws.Range("C3").Value = "maria.gonzalez+20251109_005042@test.com"
$ws.Range("C4").Value = "carlos.rodriguez+20251109_005042@test.com"
$ws.Range("C5").Value = "ana.martinez+20251109_005042@test.com"
$ws.Range("C6").Value = "luis.garcia+20251109_005042@test.com"

# "LoginData" reuses the same two e-mail addresses (Juan and Maria) in
# column A, rows 2-3; update them too so no stale shared string with the
# old timestamp is left behind.
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Range("A2").Value = "juan.perez+20251109_005042@test.com"
$wsLogin.Range("A3").Value = "maria.gonzalez+20251109_005042@test.com"
